$wb = $excel.ActiveWorkbook

# --- 1. "BAU Emissions" sheet: rename the " : NoSettings" suffix to " : test" ---
#        for every label found in column A.
$wsBau = $wb.Worksheets.Item("BAU Emissions")
$suffixOld = " : NoSettings"
$suffixNew = " : test"
$lastRow = $wsBau.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $wsBau.Cells.Item($r, 1)
    $v = $cell.Value()
    if ($v -ne $null -and $v -is [string] -and $v.EndsWith($suffixOld)) {
        $newv = $v.Substring(0, $v.Length - $suffixOld.Length) + $suffixNew
        $cell.Value = $newv
    }
}

# --- 2. Update the natural-gas / iron-and-steel emissions row (row 94), columns M:AE ---
$wsBau.Range("M94").Value = 1001080
$wsBau.Range("N94").Value = 2002150
$wsBau.Range("O94").Value = 3003230
$wsBau.Range("P94").Value = 4004300
$wsBau.Range("Q94:AE94").Value = 5005380

# --- 3. Update the "last revised" date on the About sheet ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45387

# --- 4. Re-select the ranges/cells that were active on each affected sheet ---
$wsBau.Range("A30:AE280").Select() | Out-Null

# --- 5. Make "About" the active sheet/tab (was "Current and Planned Capacity") ---
$wsAbout.Activate()
$wsAbout.Range("E29").Select() | Out-Null
